# Apply the two changes captured by the commit:
#   1. Slide 16's table switches from table style {E25C2517-19C1-4BA5-8FC7-2AF5D691195B}
#      to table style {6C38D785-7E5A-4DCF-B48F-0616FC27C1FC}.
#   2. The deck's applied theme ("Integral") has its colour scheme swapped back to
#      the stock "Office Theme" palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -------------------------------------------------
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{6C38D785-7E5A-4DCF-B48F-0616FC27C1FC}")

# --- 2. Theme colour scheme: Integral -> Office Theme ---------------------------
$colorScheme = $p.Slides.Item(1).ThemeColorScheme

function Set-ThemeColor {
    param($scheme, [int]$index, [string]$hex)
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    $scheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Index order exposed by ThemeColorScheme: 1 dk1, 2 lt1, 3 dk2, 4 lt2,
# 5-10 accent1..accent6, 11 hlink, 12 folHlink.
Set-ThemeColor $colorScheme 1  "000000"
Set-ThemeColor $colorScheme 2  "FFFFFF"
Set-ThemeColor $colorScheme 3  "44546A"
Set-ThemeColor $colorScheme 4  "E7E6E6"
Set-ThemeColor $colorScheme 5  "5B9BD5"
Set-ThemeColor $colorScheme 6  "ED7D31"
Set-ThemeColor $colorScheme 7  "A5A5A5"
Set-ThemeColor $colorScheme 8  "FFC000"
Set-ThemeColor $colorScheme 9  "4472C4"
Set-ThemeColor $colorScheme 10 "70AD47"
Set-ThemeColor $colorScheme 11 "0563C1"
Set-ThemeColor $colorScheme 12 "954F72"

# Best-effort: rename the colour scheme / design / theme back to their stock
# "Office" labels as well (harmless no-op on hosts that keep these read-only).
try { $colorScheme.Name = "Office" } catch {}
try { $p.Designs.Item(1).Name = "Office Theme" } catch {}
try { $p.SlideMaster.Theme.Name = "Office Theme" } catch {}
